$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Rename the sheet
$ws.Name = "subset_test_madeup"

# Update the defined name to point at the renamed sheet
$wb.Names.Item("WIWA_filtered").RefersTo = "=subset_test_madeup!`$A`$1:`$F`$669178"

# Update selection on the sheet
$ws.Range("D14").Select()

# Update date values in column F (rows 2-20), keeping same style (date format)
$dates = @{
    2 = 36911
    3 = 37258
    4 = 37624
    5 = 37990
    6 = 38357
    7 = 36923
    8 = 37289
    9 = 37655
    10 = 38021
    11 = 36951
    12 = 37317
    13 = 37683
    14 = 37684
    15 = 37685
    16 = 38052
    17 = 38418
    18 = 37347
    19 = 37348
    20 = 39541
}

foreach ($row in $dates.Keys) {
    $ws.Cells.Item($row, 6).Value = $dates[$row]
}
